# Apply the updated cryptocurrency market data (prices / 1h volume / row
# reordering for a few coins) produced by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price/Volume columns hold text (e.g. "1.00", "66.493.23") rather than
# numbers, so every value is written with a leading apostrophe to force Excel
# to store it as text instead of auto-converting it to a number/date.

$ws.Cells.Item(2, 4).Value = '''66.493.23'
$ws.Cells.Item(2, 5).Value = '''  +4.06%  '
$ws.Cells.Item(3, 4).Value = '''3.847.42'
$ws.Cells.Item(3, 5).Value = '''  +9.23%  '
$ws.Cells.Item(4, 5).Value = '''  +0.30%  '
$ws.Cells.Item(5, 4).Value = '''426.18'
$ws.Cells.Item(5, 5).Value = '''  +8.26%  '
$ws.Cells.Item(6, 4).Value = '''131.52'
$ws.Cells.Item(6, 5).Value = '''  +7.02%  '
$ws.Cells.Item(7, 4).Value = '''3.843.38'
$ws.Cells.Item(7, 5).Value = '''  +9.41%  '
$ws.Cells.Item(8, 4).Value = '''0.613'
$ws.Cells.Item(8, 5).Value = '''  +4.19%  '
$ws.Cells.Item(9, 4).Value = '''0.999'
$ws.Cells.Item(9, 5).Value = '''  -0.06%  '
$ws.Cells.Item(10, 4).Value = '''0.730'
$ws.Cells.Item(10, 5).Value = '''  +8.20%  '
$ws.Cells.Item(11, 4).Value = '''0.157'
$ws.Cells.Item(11, 5).Value = '''  +2.46%  '
$ws.Cells.Item(12, 5).Value = '''  -2.95%  '
$ws.Cells.Item(13, 4).Value = '''41.64'
$ws.Cells.Item(13, 5).Value = '''  +6.94%  '
$ws.Cells.Item(14, 4).Value = '''10.48'
$ws.Cells.Item(14, 5).Value = '''  +13.16%  '
$ws.Cells.Item(15, 4).Value = '''4.456.39'
$ws.Cells.Item(15, 5).Value = '''  +10.43%  '
$ws.Cells.Item(16, 4).Value = '''15.73'
$ws.Cells.Item(16, 5).Value = '''  +25.58%  '
$ws.Cells.Item(17, 4).Value = '''3.842.61'
$ws.Cells.Item(17, 5).Value = '''  +9.62%  '
$ws.Cells.Item(18, 5).Value = '''  +1.24%  '
$ws.Cells.Item(19, 4).Value = '''20.06'
$ws.Cells.Item(19, 5).Value = '''  +7.05%  '
$ws.Cells.Item(20, 4).Value = '''1.10'
$ws.Cells.Item(20, 5).Value = '''  +7.77%  '
$ws.Cells.Item(21, 4).Value = '''66.918.80'
$ws.Cells.Item(21, 5).Value = '''  +4.93%  '
$ws.Cells.Item(22, 4).Value = '''414.39'
$ws.Cells.Item(22, 5).Value = '''  +4.99%  '
$ws.Cells.Item(23, 4).Value = '''15.05'
$ws.Cells.Item(23, 5).Value = '''  +8.40%  '
$ws.Cells.Item(24, 4).Value = '''84.76'
$ws.Cells.Item(24, 5).Value = '''  +4.81%  '
$ws.Cells.Item(25, 4).Value = '''3.08'
$ws.Cells.Item(25, 5).Value = '''  +7.53%  '
$ws.Cells.Item(26, 4).Value = '''37.57'
$ws.Cells.Item(26, 5).Value = '''  +13.07%  '
$ws.Cells.Item(27, 5).Value = '''  +13.49%  '
$ws.Cells.Item(28, 5).Value = '''  +9.14%  '
$ws.Cells.Item(29, 4).Value = '''5.35'
$ws.Cells.Item(29, 5).Value = '''  +1.96%  '
$ws.Cells.Item(30, 4).Value = '''9.19'
$ws.Cells.Item(30, 5).Value = '''  +35.39%  '
$ws.Cells.Item(31, 4).Value = '''719.80'
$ws.Cells.Item(31, 5).Value = '''  +8.82%  '
$ws.Cells.Item(32, 4).Value = '''13.71'
$ws.Cells.Item(32, 5).Value = '''  +14.91%  '
$ws.Cells.Item(33, 5).Value = '''  +13.47%  '
$ws.Cells.Item(34, 5).Value = '''  +6.00%  '
$ws.Cells.Item(35, 4).Value = '''1.00'
$ws.Cells.Item(35, 5).Value = '''  -0.12%  '
$ws.Cells.Item(36, 2).Value = '''NEARProtocol'
$ws.Cells.Item(36, 3).Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(36, 4).Value = '''5.84'
$ws.Cells.Item(36, 5).Value = '''  +44.09%  '
$ws.Cells.Item(37, 2).Value = '''InjectiveProtocol'
$ws.Cells.Item(37, 3).Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(37, 4).Value = '''39.22'
$ws.Cells.Item(37, 5).Value = '''  +6.39%  '
$ws.Cells.Item(38, 2).Value = '''Kaspa'
$ws.Cells.Item(38, 3).Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(38, 4).Value = '''0.152'
$ws.Cells.Item(38, 5).Value = '''  +0.07%  '
$ws.Cells.Item(39, 4).Value = '''55.66'
$ws.Cells.Item(39, 5).Value = '''  +3.31%  '
$ws.Cells.Item(40, 4).Value = '''0.0₃0745'
$ws.Cells.Item(40, 5).Value = '''  +15.71%  '
$ws.Cells.Item(41, 4).Value = '''0.0464'
$ws.Cells.Item(41, 5).Value = '''  +5.99%  '
$ws.Cells.Item(42, 4).Value = '''2.90'
$ws.Cells.Item(42, 5).Value = '''  +7.23%  '
$ws.Cells.Item(43, 5).Value = '''  +0.69%  '
$ws.Cells.Item(44, 4).Value = '''3.27'
$ws.Cells.Item(44, 5).Value = '''  +9.19%  '
$ws.Cells.Item(45, 5).Value = '''  +3.70%  '
$ws.Cells.Item(46, 4).Value = '''3.39'
$ws.Cells.Item(46, 5).Value = '''  +10.05%  '
$ws.Cells.Item(47, 4).Value = '''0.318'
$ws.Cells.Item(47, 5).Value = '''  +14.67%  '
$ws.Cells.Item(48, 4).Value = '''141.84'
$ws.Cells.Item(48, 5).Value = '''  +1.77%  '
$ws.Cells.Item(49, 4).Value = '''2.05'
$ws.Cells.Item(49, 5).Value = '''  +4.94%  '
$ws.Cells.Item(50, 2).Value = '''Stacks'
$ws.Cells.Item(50, 3).Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(50, 4).Value = '''2.84'
$ws.Cells.Item(50, 5).Value = '''  +5.17%  '
$ws.Cells.Item(51, 2).Value = '''WEMIXToken'
$ws.Cells.Item(51, 3).Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(51, 4).Value = '''2.60'
$ws.Cells.Item(51, 5).Value = '''  +4.97%  '

# Writing a leading apostrophe makes Excel apply a "Text" quote-prefix style
# to the cell (even though none existed in the original file). Clear the
# formatting on the touched range so the cells keep their original (unstyled)
# appearance while retaining the text values just assigned.
$ws.Range("B2:E51").ClearFormats()
